$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Clear the "comentario" text for rows 5 and 6 (tasks whose comments were
# removed / cleaned up), then let Excel re-fit the row height now that the
# multi-line comments are gone.
$ws.Range("D5").ClearContents()
$ws.Range("D6").ClearContents()
$ws.Rows.Item(5).AutoFit()
$ws.Rows.Item(6).AutoFit()

# Row 14's comment got a touch-up too, and its row height was manually
# tightened to 19.5 (with an explicit custom height) rather than the old 30.
$ws.Rows.Item(14).RowHeight = 19.5

# Update the sheet view: previously scrolled to show row 4 at the top with
# D7 selected; now it shows from the top with D5 selected.
$null = $ws.Range("D5").Select()
